$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so ambiguous numeric-looking
# strings like "320.51" are not auto-coerced to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '48.062.29'
$ws.Range('E2').Value = '  +0.48%  '

$ws.Range('D3').Value = '2.499.29'
$ws.Range('E3').Value = '  +0.03%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = '320.51'
$ws.Range('E5').Value = '  -0.90%  '

$ws.Range('D6').Value = '107.46'
$ws.Range('E6').Value = '  -1.62%  '

$ws.Range('D7').Value = '0.525'
$ws.Range('E7').Value = '  +0.12%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('E9').Value = '  -2.46%  '

$ws.Range('D10').Value = '39.61'
$ws.Range('E10').Value = '  -2.80%  '

$ws.Range('D11').Value = '20.07'
$ws.Range('E11').Value = '  +7.25%  '

$ws.Range('D12').Value = '0.0811'
$ws.Range('E12').Value = '  -0.39%  '

$ws.Range('E13').Value = '  -0.20%  '

$ws.Range('D14').Value = '7.09'
$ws.Range('E14').Value = '  -1.98%  '

$ws.Range('D15').Value = '2.891.16'
$ws.Range('E15').Value = '  +0.16%  '

$ws.Range('D16').Value = '2.500.57'
$ws.Range('E16').Value = '  +0.11%  '

$ws.Range('D17').Value = '0.832'
$ws.Range('E17').Value = '  -2.44%  '

$ws.Range('D18').Value = '47.918.74'
$ws.Range('E18').Value = '  +0.43%  '

$ws.Range('D19').Value = '12.95'
$ws.Range('E19').Value = '  -2.16%  '

$ws.Range('D20').Value = '6.71'
$ws.Range('E20').Value = '  +0.93%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0939'
$ws.Range('E21').Value = '  -0.59%  '

$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').Value = '2.76'
$ws.Range('E22').Value = '  -0.95%  '

$ws.Range('E23').Value = '  +11.75%  '

$ws.Range('D24').Value = '71.41'
$ws.Range('E24').Value = '  +0.84%  '

$ws.Range('D25').Value = '2.53'
$ws.Range('E25').Value = '  -1.09%  '

$ws.Range('D27').Value = '25.88'
$ws.Range('E27').Value = '  -0.15%  '

$ws.Range('D28').Value = '9.71'
$ws.Range('E28').Value = '  -2.94%  '

$ws.Range('E29').Value = '  +1.45%  '

$ws.Range('D30').Value = '35.20'
$ws.Range('E30').Value = '  -0.47%  '

$ws.Range('D31').Value = '2.09'
$ws.Range('E31').Value = '  -4.86%  '

$ws.Range('D32').Value = '49.53'
$ws.Range('E32').Value = '  -0.55%  '

$ws.Range('D33').Value = '19.49'
$ws.Range('E33').Value = '  -2.24%  '

$ws.Range('E34').Value = '  -0.17%  '

$ws.Range('D35').Value = '5.29'
$ws.Range('E35').Value = '  -1.46%  '

$ws.Range('D36').Value = '0.0783'
$ws.Range('E36').Value = '  -1.12%  '

$ws.Range('E37').Value = '  -1.47%  '

$ws.Range('D38').Value = '4.61'
$ws.Range('E38').Value = '  -1.34%  '

$ws.Range('D39').Value = '2.87'
$ws.Range('E39').Value = '  -3.29%  '

$ws.Range('E40').Value = '  -0.80%  '

$ws.Range('D41').Value = '120.72'
$ws.Range('E41').Value = '  +1.51%  '

$ws.Range('E42').Value = '  -0.58%  '

$ws.Range('D43').Value = '21.19'
$ws.Range('E43').Value = '  -4.68%  '

$ws.Range('E44').Value = '  +1.14%  '

$ws.Range('D45').Value = '2.022.64'
$ws.Range('E45').Value = '  +0.97%  '

$ws.Range('D46').Value = '3.14'
$ws.Range('E46').Value = '  +2.22%  '

$ws.Range('E47').Value = '  -1.54%  '

$ws.Range('D48').Value = '1.83'
$ws.Range('E48').Value = '  +0.78%  '

$ws.Range('D49').Value = '8.97'
$ws.Range('E49').Value = '  -0.78%  '

$ws.Range('D50').Value = '5.17'
$ws.Range('E50').Value = '  +0.48%  '

$ws.Range('D51').Value = '80.33'
$ws.Range('E51').Value = '  +3.09%  '

# Restore the default (no explicit number format) style so the saved
# workbook matches the original "Normal"/General styling on these cells.
$ws.Range("D2:D51").Style = "Normal"
